$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ------------------------------------------------------------------
# 1. Duplicate the "en" sheet to create the new "es" sheet.
#    Copy() reproduces every row/column width, cell style and value,
#    so we only need to overwrite column B with the Spanish text below.
# ------------------------------------------------------------------
$ws1.Copy($null, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "es"

# ------------------------------------------------------------------
# 2. Translate column B (the "Value" column) on the "es" sheet.
#    Column A (the keys), C and D are left untouched.
# ------------------------------------------------------------------
$ws2.Cells.Item(2, 2).Value = 'Bienvenido!'
$ws2.Cells.Item(3, 2).Value = 'ROBERT QUEST'
$ws2.Cells.Item(4, 2).Value = 'Escapar del cuadrante'
$ws2.Cells.Item(5, 2).Value = 'CRÉDITOS'
$ws2.Cells.Item(6, 2).Value = 'Escrito por: David Dionisio\nMúsica de: Kevin Macleod'
$ws2.Cells.Item(7, 2).Value = 'JUGAR'
$ws2.Cells.Item(8, 2).Value = 'OPCIONES'
$ws2.Cells.Item(9, 2).Value = 'MÚSICA'
$ws2.Cells.Item(10, 2).Value = 'SONIDO'
$ws2.Cells.Item(11, 2).Value = 'DISCURSO'
$ws2.Cells.Item(12, 2).Value = 'ENCENDIDO'
$ws2.Cells.Item(13, 2).Value = 'APAGADO'
$ws2.Cells.Item(14, 2).Value = 'CERCA'
$ws2.Cells.Item(15, 2).Value = 'CLARO'
$ws2.Cells.Item(16, 2).Value = 'Origen'
$ws2.Cells.Item(17, 2).Value = 'El Eje X'
$ws2.Cells.Item(18, 2).Value = 'El Eje Y'
$ws2.Cells.Item(19, 2).Value = 'El Eje X (+)'
$ws2.Cells.Item(20, 2).Value = 'El Eje X (-)'
$ws2.Cells.Item(21, 2).Value = 'El Eje Y(+)'
$ws2.Cells.Item(22, 2).Value = 'El Eje Y (-)'
$ws2.Cells.Item(23, 2).Value = 'Cuadrante 1'
$ws2.Cells.Item(24, 2).Value = 'Cuadrante 2'
$ws2.Cells.Item(25, 2).Value = 'Cuadrante 3'
$ws2.Cells.Item(26, 2).Value = 'Cuadrante 4'
$ws2.Cells.Item(27, 2).Value = 'Pulsa este botón para mostrar sugerencias.'
$ws2.Cells.Item(28, 2).Value = 'Arrastre el elemento hasta el punto designado.'
$ws2.Cells.Item(29, 2).Value = 'Pulsa este botón para jugar.'
$ws2.Cells.Item(30, 2).Value = 'Para ayudar a Robert aún más, primero debemos aprender sobre el plano de coordenadas.'
$ws2.Cells.Item(31, 2).Value = 'El plano de coordenadas es una superficie bidimensional formada por dos líneas: horizontal y vertical.'
$ws2.Cells.Item(32, 2).Value = 'Primero, la línea horizontal: el eje X.'
$ws2.Cells.Item(33, 2).Value = 'A partir del origen, los valores X positivos van a la derecha.'
$ws2.Cells.Item(34, 2).Value = '... y los valores X negativos van a la izquierda.'
$ws2.Cells.Item(35, 2).Value = 'Ahora, la línea vertical: el eje Y.'
$ws2.Cells.Item(36, 2).Value = 'A partir del origen, los valores Y positivos suben.'
$ws2.Cells.Item(37, 2).Value = '... y los valores Y negativos bajan.'
$ws2.Cells.Item(38, 2).Value = 'Cuando emparejas los valores X e Y juntos, obtendás una coordenada.'
$ws2.Cells.Item(39, 2).Value = 'La coordenada es donde se cruzan las líneas de los ejes X e Y. Observa que los dos valores cambian a medida que Robert se mueve.'
$ws2.Cells.Item(40, 2).Value = 'A continuación, repasaremos los cuadrantes.'
$ws2.Cells.Item(41, 2).Value = 'Los cuadrantes son las cuatro secciones del plano de coordenadas.'
$ws2.Cells.Item(42, 2).Value = 'Como puedes ver, los cuadrantes están divididos por los ejes X e Y desde el origen.'
$ws2.Cells.Item(43, 2).Value = 'Cada uno determina los valores de signo del X e Y: positivos o negativos.'
$ws2.Cells.Item(44, 2).Value = 'Ahora vamos a ayudar a Robert a pasar del cuadrante 1 al cuadrante 2.'
$ws2.Cells.Item(45, 2).Value = '¡Genial! El próximo paso: Cuadrante 3.'
$ws2.Cells.Item(46, 2).Value = 'Un último destino: Cuadrante 4.'
$ws2.Cells.Item(47, 2).Value = '!Excelente! Ahora estás listo para guiar a Robert de vuelta a su familia.'
$ws2.Cells.Item(48, 2).Value = 'Ahora vamos a repasar la reflexión.'
$ws2.Cells.Item(49, 2).Value = 'Un punto se refleja al voltear los signos de sus valores: positivo a negativo y viceversa.'
$ws2.Cells.Item(50, 2).Value = 'Aquí puedes ver el valor X reflejado en el eje Y.'
$ws2.Cells.Item(51, 2).Value = '... y el valor Y, que se refleja en el eje X.'
$ws2.Cells.Item(52, 2).Value = 'Al reflejar los valores X e Y, puede ver que la línea que conecta ambos puntos recorre el origen.'
$ws2.Cells.Item(53, 2).Value = 'Los dos puntos tienen la misma distancia del origen.'
$ws2.Cells.Item(54, 2).Value = '¡Ahora pongamos esto en práctica! Coloque los reflejos coincidentes en las fichas resaltadas.'
$ws2.Cells.Item(55, 2).Value = 'Este es Robert. Ha estado corriendo en círculos durante años.'
$ws2.Cells.Item(56, 2).Value = 'Por favor, rompe este ciclo para que pueda reunirse con su familia.'
$ws2.Cells.Item(57, 2).Value = '¡Excelente!  Con tu dirección, Robert finalmente puede seguir adelante.'
$ws2.Cells.Item(58, 2).Value = 'Después de un largo y arduo viaje, Robert finalmente se une a su familia.'
$ws2.Cells.Item(59, 2).Value = '¡Gracias por jugar!'
$ws2.Cells.Item(60, 2).Value = 'COMPLETO'
$ws2.Cells.Item(61, 2).Value = 'NIVEL COMPLETO'
$ws2.Cells.Item(62, 2).Value = 'CUADRANTE'
$ws2.Cells.Item(63, 2).Value = '¿En qué cuadrante se teletransportará Robert?'

# ------------------------------------------------------------------
# 3. Restore view/selection state for both sheets.
#    "en" keeps being the active/selected tab (Copy() would otherwise
#    leave the new sheet active), scrolled back to the top with B21
#    selected; "es" opens with B6 selected.
# ------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("B21").Select()
$ws2.Range("B6").Select()
$ws1.Activate()

# ------------------------------------------------------------------
# 4. Best-effort restore of the application window geometry.
# ------------------------------------------------------------------
$w = $excel.ActiveWindow
$w.Left = 4785
$w.Top = 1830
$w.Width = 28365
$w.Height = 18045
